$d = $word.ActiveDocument

$d.Content.Find.Execute("67+31=98", $true, $false, $false, $false, $false, $true, 1, $false, "99-33=66", 2) | Out-Null
$d.Content.Find.Execute("54-23=31", $true, $false, $false, $false, $false, $true, 1, $false, "5+33=38", 2) | Out-Null
$d.Content.Find.Execute("11+22=33", $true, $false, $false, $false, $false, $true, 1, $false, "25-2=23", 2) | Out-Null
$d.Content.Find.Execute("97+1=98", $true, $false, $false, $false, $false, $true, 1, $false, "89-68=21", 2) | Out-Null
$d.Content.Find.Execute("15+78=93", $true, $false, $false, $false, $false, $true, 1, $false, "71-69=2", 2) | Out-Null
$d.Content.Find.Execute("82-60=22", $true, $false, $false, $false, $false, $true, 1, $false, "83-40=43", 2) | Out-Null
$d.Content.Find.Execute("33+65=98", $true, $false, $false, $false, $false, $true, 1, $false, "7+54=61", 2) | Out-Null
$d.Content.Find.Execute("58-5=53", $true, $false, $false, $false, $false, $true, 1, $false, "71-43=28", 2) | Out-Null
$d.Content.Find.Execute("3+43=46", $true, $false, $false, $false, $false, $true, 1, $false, "97-81=16", 2) | Out-Null
$d.Content.Find.Execute("38+49=87", $true, $false, $false, $false, $false, $true, 1, $false, "8+28=36", 2) | Out-Null
$d.Content.Find.Execute("36-16=20", $true, $false, $false, $false, $false, $true, 1, $false, "36+53=89", 2) | Out-Null
$d.Content.Find.Execute("78-76=2", $true, $false, $false, $false, $false, $true, 1, $false, "46-45=1", 2) | Out-Null
$d.Content.Find.Execute("75+3=78", $true, $false, $false, $false, $false, $true, 1, $false, "23+27=50", 2) | Out-Null
$d.Content.Find.Execute("2+68=70", $true, $false, $false, $false, $false, $true, 1, $false, "50-29=21", 2) | Out-Null
$d.Content.Find.Execute("73-1=72", $true, $false, $false, $false, $false, $true, 1, $false, "27+49=76", 2) | Out-Null
$d.Content.Find.Execute("62+28=90", $true, $false, $false, $false, $false, $true, 1, $false, "0+22=22", 2) | Out-Null
$d.Content.Find.Execute("29-17=12", $true, $false, $false, $false, $false, $true, 1, $false, "57-45=12", 2) | Out-Null
$d.Content.Find.Execute("44+47=91", $true, $false, $false, $false, $false, $true, 1, $false, "14+80=94", 2) | Out-Null
$d.Content.Find.Execute("72+16=88", $true, $false, $false, $false, $false, $true, 1, $false, "56+7=63", 2) | Out-Null
$d.Content.Find.Execute("82-65=17", $true, $false, $false, $false, $false, $true, 1, $false, "27-2=25", 2) | Out-Null
$d.Content.Find.Execute("40+16=56", $true, $false, $false, $false, $false, $true, 1, $false, "84-56=28", 2) | Out-Null
$d.Content.Find.Execute("29+70=99", $true, $false, $false, $false, $false, $true, 1, $false, "75-43=32", 2) | Out-Null
$d.Content.Find.Execute("4+4=8", $true, $false, $false, $false, $false, $true, 1, $false, "0+39=39", 2) | Out-Null
$d.Content.Find.Execute("89-19=70", $true, $false, $false, $false, $false, $true, 1, $false, "96-61=35", 2) | Out-Null
$d.Content.Find.Execute("66-45=21", $true, $false, $false, $false, $false, $true, 1, $false, "90-29=61", 2) | Out-Null
$d.Content.Find.Execute("96-69=27", $true, $false, $false, $false, $false, $true, 1, $false, "11+13=24", 2) | Out-Null
$d.Content.Find.Execute("39+28=67", $true, $false, $false, $false, $false, $true, 1, $false, "50+34=84", 2) | Out-Null
$d.Content.Find.Execute("63+30=93", $true, $false, $false, $false, $false, $true, 1, $false, "82-73=9", 2) | Out-Null
$d.Content.Find.Execute("97-76=21", $true, $false, $false, $false, $false, $true, 1, $false, "57+20=77", 2) | Out-Null
$d.Content.Find.Execute("91-62=29", $true, $false, $false, $false, $false, $true, 1, $false, "55+24=79", 2) | Out-Null
$d.Content.Find.Execute("39-35=4", $true, $false, $false, $false, $false, $true, 1, $false, "29+14=43", 2) | Out-Null
$d.Content.Find.Execute("77-51=26", $true, $false, $false, $false, $false, $true, 1, $false, "27-22=5", 2) | Out-Null
$d.Content.Find.Execute("15+4=19", $true, $false, $false, $false, $false, $true, 1, $false, "35+62=97", 2) | Out-Null
$d.Content.Find.Execute("11+6=17", $true, $false, $false, $false, $false, $true, 1, $false, "30-17=13", 2) | Out-Null
$d.Content.Find.Execute("30+32=62", $true, $false, $false, $false, $false, $true, 1, $false, "55+29=84", 2) | Out-Null
$d.Content.Find.Execute("5-2=3", $true, $false, $false, $false, $false, $true, 1, $false, "11+20=31", 2) | Out-Null
$d.Content.Find.Execute("30+62=92", $true, $false, $false, $false, $false, $true, 1, $false, "49+44=93", 2) | Out-Null
$d.Content.Find.Execute("95-31=64", $true, $false, $false, $false, $false, $true, 1, $false, "80+10=90", 2) | Out-Null
$d.Content.Find.Execute("90-74=16", $true, $false, $false, $false, $false, $true, 1, $false, "66-57=9", 2) | Out-Null
$d.Content.Find.Execute("54-25=29", $true, $false, $false, $false, $false, $true, 1, $false, "26+6=32", 2) | Out-Null
$d.Content.Find.Execute("77-14=63", $true, $false, $false, $false, $false, $true, 1, $false, "54+26=80", 2) | Out-Null
$d.Content.Find.Execute("42-28=14", $true, $false, $false, $false, $false, $true, 1, $false, "22+66=88", 2) | Out-Null
$d.Content.Find.Execute("61-44=17", $true, $false, $false, $false, $false, $true, 1, $false, "78-66=12", 2) | Out-Null
$d.Content.Find.Execute("66-9=57", $true, $false, $false, $false, $false, $true, 1, $false, "96-9=87", 2) | Out-Null
$d.Content.Find.Execute("2+83=85", $true, $false, $false, $false, $false, $true, 1, $false, "47-5=42", 2) | Out-Null
$d.Content.Find.Execute("98-79=19", $true, $false, $false, $false, $false, $true, 1, $false, "45-28=17", 2) | Out-Null
$d.Content.Find.Execute("18+2=20", $true, $false, $false, $false, $false, $true, 1, $false, "16+76=92", 2) | Out-Null
$d.Content.Find.Execute("34+33=67", $true, $false, $false, $false, $false, $true, 1, $false, "33+44=77", 2) | Out-Null
$d.Content.Find.Execute("43+26=69", $true, $false, $false, $false, $false, $true, 1, $false, "0+76=76", 2) | Out-Null
$d.Content.Find.Execute("79-5=74", $true, $false, $false, $false, $false, $true, 1, $false, "53-41=12", 2) | Out-Null
$d.Content.Find.Execute("98-61=37", $true, $false, $false, $false, $false, $true, 1, $false, "85-45=40", 2) | Out-Null
$d.Content.Find.Execute("52-15=37", $true, $false, $false, $false, $false, $true, 1, $false, "56+20=76", 2) | Out-Null
$d.Content.Find.Execute("0+4=4", $true, $false, $false, $false, $false, $true, 1, $false, "0+75=75", 2) | Out-Null
$d.Content.Find.Execute("91-85=6", $true, $false, $false, $false, $false, $true, 1, $false, "46-1=45", 2) | Out-Null
$d.Content.Find.Execute("86-30=56", $true, $false, $false, $false, $false, $true, 1, $false, "75+5=80", 2) | Out-Null
$d.Content.Find.Execute("29-7=22", $true, $false, $false, $false, $false, $true, 1, $false, "28+43=71", 2) | Out-Null
$d.Content.Find.Execute("33+39=72", $true, $false, $false, $false, $false, $true, 1, $false, "14+24=38", 2) | Out-Null
$d.Content.Find.Execute("28-23=5", $true, $false, $false, $false, $false, $true, 1, $false, "34-17=17", 2) | Out-Null
$d.Content.Find.Execute("83-67=16", $true, $false, $false, $false, $false, $true, 1, $false, "44+0=44", 2) | Out-Null
$d.Content.Find.Execute("84+14=98", $true, $false, $false, $false, $false, $true, 1, $false, "42+31=73", 2) | Out-Null
$d.Content.Find.Execute("5+46=51", $true, $false, $false, $false, $false, $true, 1, $false, "92-2=90", 2) | Out-Null
$d.Content.Find.Execute("43-19=24", $true, $false, $false, $false, $false, $true, 1, $false, "54+33=87", 2) | Out-Null
$d.Content.Find.Execute("69-17=52", $true, $false, $false, $false, $false, $true, 1, $false, "51-28=23", 2) | Out-Null
$d.Content.Find.Execute("82-78=4", $true, $false, $false, $false, $false, $true, 1, $false, "27+66=93", 2) | Out-Null
$d.Content.Find.Execute("86+10=96", $true, $false, $false, $false, $false, $true, 1, $false, "13+45=58", 2) | Out-Null
$d.Content.Find.Execute("83-37=46", $true, $false, $false, $false, $false, $true, 1, $false, "19-5=14", 2) | Out-Null
$d.Content.Find.Execute("87+3=90", $true, $false, $false, $false, $false, $true, 1, $false, "5+90=95", 2) | Out-Null
$d.Content.Find.Execute("7+69=76", $true, $false, $false, $false, $false, $true, 1, $false, "55-21=34", 2) | Out-Null
$d.Content.Find.Execute("77-24=53", $true, $false, $false, $false, $false, $true, 1, $false, "98-89=9", 2) | Out-Null
$d.Content.Find.Execute("8+17=25", $true, $false, $false, $false, $false, $true, 1, $false, "11-5=6", 2) | Out-Null
$d.Content.Find.Execute("98-36=62", $true, $false, $false, $false, $false, $true, 1, $false, "97-16=81", 2) | Out-Null
$d.Content.Find.Execute("74-17=57", $true, $false, $false, $false, $false, $true, 1, $false, "45+37=82", 2) | Out-Null
$d.Content.Find.Execute("46+12=58", $true, $false, $false, $false, $false, $true, 1, $false, "12+2=14", 2) | Out-Null
$d.Content.Find.Execute("35+12=47", $true, $false, $false, $false, $false, $true, 1, $false, "26+38=64", 2) | Out-Null
$d.Content.Find.Execute("60+12=72", $true, $false, $false, $false, $false, $true, 1, $false, "52-44=8", 2) | Out-Null
$d.Content.Find.Execute("9+64=73", $true, $false, $false, $false, $false, $true, 1, $false, "14+85=99", 2) | Out-Null
$d.Content.Find.Execute("55-48=7", $true, $false, $false, $false, $false, $true, 1, $false, "88-41=47", 2) | Out-Null
$d.Content.Find.Execute("80-2=78", $true, $false, $false, $false, $false, $true, 1, $false, "4-2=2", 2) | Out-Null
$d.Content.Find.Execute("84-84=0", $true, $false, $false, $false, $false, $true, 1, $false, "68-27=41", 2) | Out-Null
$d.Content.Find.Execute("14+12=26", $true, $false, $false, $false, $false, $true, 1, $false, "82-63=19", 2) | Out-Null
$d.Content.Find.Execute("25+67=92", $true, $false, $false, $false, $false, $true, 1, $false, "72-12=60", 2) | Out-Null
$d.Content.Find.Execute("76-28=48", $true, $false, $false, $false, $false, $true, 1, $false, "4+9=13", 2) | Out-Null
$d.Content.Find.Execute("63-21=42", $true, $false, $false, $false, $false, $true, 1, $false, "6+1=7", 2) | Out-Null
$d.Content.Find.Execute("79-78=1", $true, $false, $false, $false, $false, $true, 1, $false, "63+11=74", 2) | Out-Null
$d.Content.Find.Execute("10+86=96", $true, $false, $false, $false, $false, $true, 1, $false, "57+14=71", 2) | Out-Null
$d.Content.Find.Execute("83-23=60", $true, $false, $false, $false, $false, $true, 1, $false, "5-0=5", 2) | Out-Null
$d.Content.Find.Execute("40-2=38", $true, $false, $false, $false, $false, $true, 1, $false, "0+34=34", 2) | Out-Null
$d.Content.Find.Execute("99-31=68", $true, $false, $false, $false, $false, $true, 1, $false, "54+10=64", 2) | Out-Null
$d.Content.Find.Execute("54-41=13", $true, $false, $false, $false, $false, $true, 1, $false, "11+81=92", 2) | Out-Null
$d.Content.Find.Execute("78-8=70", $true, $false, $false, $false, $false, $true, 1, $false, "83-59=24", 2) | Out-Null
$d.Content.Find.Execute("53+17=70", $true, $false, $false, $false, $false, $true, 1, $false, "30-16=14", 2) | Out-Null
$d.Content.Find.Execute("67-43=24", $true, $false, $false, $false, $false, $true, 1, $false, "24+47=71", 2) | Out-Null
$d.Content.Find.Execute("5+7=12", $true, $false, $false, $false, $false, $true, 1, $false, "12+72=84", 2) | Out-Null
$d.Content.Find.Execute("7+56=63", $true, $false, $false, $false, $false, $true, 1, $false, "82+2=84", 2) | Out-Null
$d.Content.Find.Execute("46-22=24", $true, $false, $false, $false, $false, $true, 1, $false, "47+9=56", 2) | Out-Null
$d.Content.Find.Execute("13+24=37", $true, $false, $false, $false, $false, $true, 1, $false, "44+32=76", 2) | Out-Null
$d.Content.Find.Execute("40-8=32", $true, $false, $false, $false, $false, $true, 1, $false, "36-24=12", 2) | Out-Null
$d.Content.Find.Execute("69-59=10", $true, $false, $false, $false, $false, $true, 1, $false, "85-19=66", 2) | Out-Null
$d.Content.Find.Execute("56-13=43", $true, $false, $false, $false, $false, $true, 1, $false, "96-87=9", 2) | Out-Null
$d.Content.Find.Execute("74+3=77", $true, $false, $false, $false, $false, $true, 1, $false, "34-15=19", 2) | Out-Null
